# Cosmetic changes to reports.
#
# 1) Rename the worksheet from "Sheet1" to "Requests".
# 2) Refresh the "Stock Control System Requests" table:
#    - Drop the "High / Dispatch note changes..." request and replace it
#      with a new "Medium / Proforma customer alert" request.
#    - Replace the "Low / Make cosmetic changes to reports" request (row 6)
#      with a new "Medium / Products PO stock split per order + date"
#      request.
#    - The "Low / Make cosmetic changes to reports" request moves down into
#      what used to be the first blank row, with a refreshed completion
#      date and an expanded note.
#    - Drop 2 now-unused trailing blank rows so the table keeps its size.
#    - Switch the "Estimated Completion Date" column to a long-date
#      display format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- 1. Rename the sheet -------------------------------------------------
$ws.Name = "Requests"

# --- 2. Turn the first blank filler row (row 7) into the "Low" request --
# Give it the border/fill look that the "Low" row (row 6) currently has.
$ws.Range("B6:E6").Copy()
$ws.Range("B7:E7").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("B7").Value2 = "Low"
$ws.Range("C7").Value2 = "Make cosmetic changes to reports "
$ws.Range("D7").Value2 = 42660
$ws.Range("E7").Value2 = "Use a template excel sheet for each report"

# --- 3. Update row 5 : High/Dispatch... -> Medium/Proforma customer alert
$ws.Range("B5").Value2 = "Medium"
$ws.Range("C5").Value2 = "Proforma customer alert"
$ws.Range("D5").Value2 = 42650
$ws.Range("E5").Value2 = "Alert when selecting customer for a new order"

# --- 4. Update row 6 : Low/Make cosmetic... -> Medium/Products PO split -
$ws.Range("B6").Value2 = "Medium"
$ws.Range("C6").Value2 = "Products PO stock split per order + date"
$ws.Range("D6").Value2 = 42653
$ws.Range("E6").ClearContents()

# --- 5. Drop 2 now-unused trailing blank rows ----------------------------
$ws.Range("B28:B29").EntireRow.Delete()

# --- 6. Long-date display for the "Estimated Completion Date" column ----
$ws.Range("D4:D30").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"

$ws.Range("B1").Select()
